# Hotfix: Mon Oct 21 17:53:26 RTZ 2024
#
# - SQL sheet: remove the stray test row 13 (hgfhgfhgfh... placeholder data)
# - Python sheet: remove the stray test row 38 (hgfhgfhgfh... placeholder data)
# - Links sheet: remove the stray test row 1 (add_links_command... placeholder data)
# - Bash sheet: replace the two placeholder rows (71/72) with real "mount"/"tar"/"mv"
#   content and add two more real rows (73/74) for the split tar/mv commands

$wb = $excel.ActiveWorkbook

# ---- SQL sheet: delete the placeholder row (row 13) ----
$wsSql = $wb.Worksheets.Item("SQL")
$wsSql.Rows.Item(13).Delete()

# ---- Python sheet: delete the placeholder row (row 38) ----
$wsPython = $wb.Worksheets.Item("Python")
$wsPython.Rows.Item(38).Delete()

# ---- Links sheet: delete the placeholder row (row 1) ----
$wsLinks = $wb.Worksheets.Item("Links")
$wsLinks.Rows.Item(1).Delete()

# ---- Bash sheet: fix up rows 71-72 and append rows 73-74 ----
$wsBash = $wb.Worksheets.Item("Bash")

$wsBash.Cells.Item(71, 1).Value = 112
$wsBash.Cells.Item(71, 2).Value = "mount"
$wsBash.Cells.Item(71, 3).Value = "Вывод всех подключенных дисков"

$wsBash.Cells.Item(72, 1).Value = 113
$wsBash.Cells.Item(72, 2).Value = "tar -cf flask-project.tar * && mv flask-project.tar /o"
$wsBash.Cells.Item(72, 3).Value = "Скрипт архивации файлов и перенос архива на другой диск"

$wsBash.Cells.Item(73, 1).Value = 114
$wsBash.Cells.Item(73, 2).Value = "tar -cf flask-project.tar *"
$wsBash.Cells.Item(73, 3).Value = "Создание архива всех файлов в папке"

$wsBash.Cells.Item(74, 1).Value = 115
$wsBash.Cells.Item(74, 2).Value = "mv flask-project.tar /o"
$wsBash.Cells.Item(74, 3).Value = "Перенос файла в другой диск"
